# fix typos in sample arc (MSEval_Heat isa.assay.xlsx)
#
# - MSEval!A4 was a duplicated "C2_measured" -> should be "C3_measured"
# - MSEval!K2:K7 ("Sample Name") pointed at "sampleOut.txt" -> should be
#   "sampleOutHeat.txt"
# - the B/C/D/H/I/J helper columns of the MSEval table carried stray
#   (empty, styled) cells left over from the template; they are cleared
#   out entirely, and the formatting on the remaining E/F/G annotation
#   cells is dropped along with them
# - MSEval's selection moves off the old E20 leftover onto the Sample
#   Name column that was just retyped

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSEval")

# Row 4 duplicated "C2_measured" -> fix to "C3_measured"
$ws.Range("A4").Value = "C3_measured"

# Sample Name column: sampleOut.txt -> sampleOutHeat.txt
$ws.Range("K2:K7").Value = "sampleOutHeat.txt"

# Drop the leftover empty styled cells in the hidden helper columns
$ws.Range("B2:D7").Clear()
$ws.Range("H2:J7").Clear()

# The remaining annotation cells (E:G) keep their text but lose the
# inherited formatting that used to be shared with the cleared columns
$ws.Range("E2:G7").ClearFormats()

# Reflect the edit in the active selection
$ws.Range("K4:K7").Select() | Out-Null
